$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.483.42"
$ws.Range("E2").Value = "  +0.20%  "
$ws.Range("D3").Value = "2.640.37"
$ws.Range("E3").Value = "  +1.23%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "536.88"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.08%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "144.82"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.07%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.998"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.571"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.62%  "
$ws.Range("D9").Value = "2.654.05"
$ws.Range("E9").Value = "  +1.45%  "
$ws.Range("E10").Value = "  +2.83%  "
$ws.Range("E11").Value = "  -0.90%  "
$ws.Range("E12").Value = "  +0.23%  "
$ws.Range("E13").Value = "  -0.27%  "
$ws.Range("D14").Value = "3.105.81"
$ws.Range("E14").Value = "  +1.09%  "
$ws.Range("D15").Value = "59.391.03"
$ws.Range("E15").Value = "  +0.18%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "21.14"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.88%  "
$ws.Range("D17").Value = "2.651.38"
$ws.Range("E17").Value = "  -0.18%  "
$ws.Range("E18").Value = "  +0.53%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "339.18"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.14%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.39"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.73%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.36"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.56%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.31"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.21%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "67.05"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.32%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.416"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.87%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.165"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.95%  "
$ws.Range("E27").Value = "  +0.02%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.27"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.52%  "
$ws.Range("D29").Value = "0.0₃0744"
$ws.Range("E29").Value = "  +0.20%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.999"
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.65"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.73%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.84"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.27%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "18.87"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.08%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "151.44"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.81%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.99"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.35%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.14"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.47%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.838"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.63%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.833"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.24%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "287.81"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +4.66%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.59"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.05%  "
$ws.Range("E42").Value = "  -0.04%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.604"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.18%  "
$ws.Range("E44").Value = "  +0.01%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "19.32"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +3.72%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0536"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.39%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0946"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.68%  "
$ws.Range("D48").Value = "1.967.95"
$ws.Range("E48").Value = "  +0.86%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0226"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.15%  "
$ws.Range("E50").Value = "  +1.28%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "18.27"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.25%  "
